$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bus")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("J2") "-19.908076343375843"
Set-TextValue $ws.Range("K2") "-57.80509971944206"
Set-TextValue $ws.Range("J3") "19.908076343375843"
Set-TextValue $ws.Range("K3") "57.80509971944206"
